$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===================================================================
# Data changes: a new "disp" scenario column (D) and a new "disp"
# scenario column (H) are populated (they were previously blank,
# duplicating the scenario values already present in columns E/I),
# plus a handful of value corrections, plus the "BEWARE" note moves
# from Q1 into L1 now that the trailing blank columns are gone.
# ===================================================================

$ws.Range("L1").Value = " BEWARE that if add a row here also need to add to set_run_inputs.r "
$ws.Range("Q1").ClearContents()

$ws.Range("D2").Value = 0.1
$ws.Range("H2").Value = 0.1

$ws.Range("D3").Value = 0.9
$ws.Range("H3").Value = 0.9

$ws.Range("D4").Value = 0.01
$ws.Range("E4").Value = 0.01
$ws.Range("F4").Value = 0.0
$ws.Range("H4").Value = 0.01
$ws.Range("I4").Value = 0.01

$ws.Range("D5").Value = 0.3
$ws.Range("E5").Value = 0.3
$ws.Range("F5").Value = 0.0
$ws.Range("H5").Value = 0.3
$ws.Range("I5").Value = 0.3

$ws.Range("D6").Value = 500.0
$ws.Range("H6").Value = 500.0

$ws.Range("D7").Value = 10.0
$ws.Range("H7").Value = 10.0

$ws.Range("D8").Value = 10.0
$ws.Range("H8").Value = 10.0

$ws.Range("D9").Value = 0.0
$ws.Range("H9").Value = 0.0

$ws.Range("D10").Value = 0.0
$ws.Range("H10").Value = 0.0

$ws.Range("D11").Value = 0.0
$ws.Range("H11").Value = 0.0

$ws.Range("D12").Value = "freq"
$ws.Range("H12").Value = "freq"

$ws.Range("D13").Value = 0.9
$ws.Range("H13").Value = 0.9

$ws.Range("D14").Value = 0.5
$ws.Range("H14").Value = 0.5

$ws.Range("D15").Value = 1.0
$ws.Range("E15").Value = 1.0
$ws.Range("H15").Value = 1.0
$ws.Range("I15").Value = 1.0

$ws.Range("D16").Value = 0.0
$ws.Range("H16").Value = 0.0

$ws.Range("D17").Value = 0.0
$ws.Range("F17").Value = 0.0
$ws.Range("H17").Value = 0.0

$ws.Range("D18").Value = 0.4
$ws.Range("H18").Value = 0.4

$ws.Range("D19").Value = 0.9
$ws.Range("H19").Value = 0.9

$ws.Range("B20").Value = 0.0
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.0
$ws.Range("E20").Value = 0.0
$ws.Range("F20").Value = 0.0
$ws.Range("G20").Value = 0.0
$ws.Range("H20").Value = 0.0
$ws.Range("I20").Value = 0.0

$ws.Range("D21").Value = 1.0
$ws.Range("H21").Value = 1.0

$ws.Range("D22").Value = 0.7
$ws.Range("E22").Value = 0.7
$ws.Range("H22").Value = 0.7

$ws.Range("D23").Value = 1.0
$ws.Range("H23").Value = 1.0

$ws.Range("D24").Value = 0.1
$ws.Range("H24").Value = 0.1

$ws.Range("D25").Value = 0.9
$ws.Range("H25").Value = 0.9

$ws.Range("D26").Value = 0.0
$ws.Range("H26").Value = 0.0

$ws.Range("D27").Value = 1.0
$ws.Range("H27").Value = 1.0

$ws.Range("D28").Value = 0.0
$ws.Range("H28").Value = 0.0

$ws.Range("D29").Value = 1.0
$ws.Range("H29").Value = 1.0

$ws.Range("D30").Value = 0.001
$ws.Range("F30").Value = 0.001
$ws.Range("H30").Value = 0.001

$ws.Range("D31").Value = 0.1
$ws.Range("H31").Value = 0.1

$ws.Range("D32").Value = 2.0
$ws.Range("E32").Value = 2.0
$ws.Range("H32").Value = 2.0
$ws.Range("I32").Value = 2.0

$ws.Range("D33").Value = 5.0
$ws.Range("H33").Value = 5.0

$ws.Range("D34").Value = 10.0
$ws.Range("H34").Value = 10.0

$ws.Range("D35").Value = 10.0
$ws.Range("H35").Value = 10.0

$ws.Range("D36").Value = 0.0

$ws.Range("D37").Value = 2.0
$ws.Range("H37").Value = 2.0

# ===================================================================
# Column widths (cosmetic resize of B,C,D,E,G,H,I to make room for
# the newly-populated D/H columns). Column.ColumnWidth is offset
# from the stored worksheet width by 5/6 in this runtime.
# ===================================================================

$ws.Columns.Item(2).ColumnWidth = 25.498697916666668
$ws.Columns.Item(3).ColumnWidth = 23.299479166666668
$ws.Columns.Item(4).ColumnWidth = 24.096354166666668
$ws.Columns.Item(5).ColumnWidth = 21.565104166666668
$ws.Columns.Item(7).ColumnWidth = 23.697916666666668
$ws.Columns.Item(8).ColumnWidth = 24.365885416666668
$ws.Columns.Item(9).ColumnWidth = 24.096354166666668

# ===================================================================
# View: freeze the header row only (was previously split at row 20 /
# col 4), and leave the selection on F36 as in the saved file.
# ===================================================================

$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("F36").Select()
